# Inserta 2 filas nuevas en la posición 775 (desplazando las filas existentes
# 775-840 hacia abajo, a 777-842) y completa los datos de la actualización
# semanal para Zanahoria, Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(775).EntireRow.Insert()
$ws.Rows.Item(775).EntireRow.Insert()

# Fila 775: Zanahoria, Primera, Región Metropolitana
$ws.Cells.Item(775, 1).Value = 6
$ws.Cells.Item(775, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(775, 3).Value = "Metropolitana"
$ws.Cells.Item(775, 4).Value = 44461
$ws.Cells.Item(775, 5).Value = 13
$ws.Cells.Item(775, 6).Value = 100114013
$ws.Cells.Item(775, 7).Value = "Zanahoria"
$ws.Cells.Item(775, 8).Value = "Sin especificar"
$ws.Cells.Item(775, 9).Value = "Primera"
$ws.Cells.Item(775, 10).Value = 2100
$ws.Cells.Item(775, 11).Value = 5000
$ws.Cells.Item(775, 12).Value = 5500
$ws.Cells.Item(775, 13).Value = 5155
$ws.Cells.Item(775, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(775, 15).Value = "Región Metropolitana"
$ws.Cells.Item(775, 16).Value = 258
$ws.Cells.Item(775, 17).Value = 20
$ws.Cells.Item(775, 18).Value = "Hortaliza"

# Fila 776: Zanahoria, Segunda, Región Metropolitana
$ws.Cells.Item(776, 1).Value = 6
$ws.Cells.Item(776, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(776, 3).Value = "Metropolitana"
$ws.Cells.Item(776, 4).Value = 44461
$ws.Cells.Item(776, 5).Value = 13
$ws.Cells.Item(776, 6).Value = 100114013
$ws.Cells.Item(776, 7).Value = "Zanahoria"
$ws.Cells.Item(776, 8).Value = "Sin especificar"
$ws.Cells.Item(776, 9).Value = "Segunda"
$ws.Cells.Item(776, 10).Value = 680
$ws.Cells.Item(776, 11).Value = 4000
$ws.Cells.Item(776, 12).Value = 4000
$ws.Cells.Item(776, 13).Value = 4000
$ws.Cells.Item(776, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(776, 15).Value = "Región Metropolitana"
$ws.Cells.Item(776, 16).Value = 200
$ws.Cells.Item(776, 17).Value = 20
$ws.Cells.Item(776, 18).Value = "Hortaliza"
